$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17 (keep header row 1 and first data row 2)
$ws.Range("A3:B17").EntireRow.Delete()

# Update row 2 values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 1.111105389870137
